$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 21539936.0074994
$ws.Range("D2").Value = 23479964.5573357
$ws.Range("E2").Value = 22808452.4929323
$ws.Range("F2").Value = 20271419.5220664
$ws.Range("G2").Value = 19599907.4576631
$ws.Range("C3").Value = 20413770.6013595
$ws.Range("D3").Value = 22357792.7025973
$ws.Range("E3").Value = 21684898.3296985
$ws.Range("F3").Value = 19142642.8730206
$ws.Range("G3").Value = 18469748.5001218
$ws.Range("C4").Value = 24325953.0976278
$ws.Range("D4").Value = 26761721.2135591
$ws.Range("E4").Value = 25918616.2623073
$ws.Range("F4").Value = 22733289.9329483
$ws.Range("G4").Value = 21890184.9816965
$ws.Range("C5").Value = 22993466.3485849
$ws.Range("D5").Value = 25403233.8261656
$ws.Range("E5").Value = 24569128.6096528
$ws.Range("F5").Value = 21417804.087517
$ws.Range("G5").Value = 20583698.8710042
$ws.Range("C6").Value = 26691951.4191559
$ws.Range("D6").Value = 29608731.6736688
$ws.Range("E6").Value = 28599131.515834
$ws.Range("F6").Value = 24784771.3224779
$ws.Range("G6").Value = 23775171.164643
$ws.Range("C7").Value = 26989964.0105518
$ws.Range("D7").Value = 30055322.4976857
$ws.Range("E7").Value = 28994294.191682
$ws.Range("F7").Value = 24985633.8294216
$ws.Range("G7").Value = 23924605.523418
$ws.Range("C8").Value = 26948630.7647638
$ws.Range("D8").Value = 30120930.2901846
$ws.Range("E8").Value = 29022885.9323318
$ws.Range("F8").Value = 24874375.5971958
$ws.Range("G8").Value = 23776331.2393431
$ws.Range("C9").Value = 24091579.3491059
$ws.Range("D9").Value = 27023985.6473797
$ws.Range("E9").Value = 26008976.7666138
$ws.Range("F9").Value = 22174181.931598
$ws.Range("G9").Value = 21159173.0508322
$ws.Range("C10").Value = 20523492.4086428
$ws.Range("D10").Value = 23101144.3982264
$ws.Range("E10").Value = 22208928.4517215
$ws.Range("F10").Value = 18838056.3655641
$ws.Range("G10").Value = 17945840.4190592
$ws.Range("C11").Value = 20011748.6685998
$ws.Range("D11").Value = 22600389.9552539
$ws.Range("E11").Value = 21704370.2268084
$ws.Range("F11").Value = 18319127.1103912
$ws.Range("G11").Value = 17423107.3819458
$ws.Range("C12").Value = 21177435.4858385
$ws.Range("D12").Value = 23994279.1915141
$ws.Range("E12").Value = 23019270.5855534
$ws.Range("F12").Value = 19335600.3861236
$ws.Range("G12").Value = 18360591.7801629
$ws.Range("C13").Value = 20855799.1096099
$ws.Range("D13").Value = 23704077.7781743
$ws.Range("E13").Value = 22718188.4267598
$ws.Range("F13").Value = 18993409.79246
$ws.Range("G13").Value = 18007520.4410455
